# Before this edit the deck ships two theme parts:
#   ppt/theme/theme1.xml -> clrScheme "Office"      (a:theme name="Office Theme")
#   ppt/theme/theme2.xml -> clrScheme "Red Violet"   (a:theme name="Integral")
#
# theme2.xml is the theme actually wired to the (only) slide master / the
# presentation, i.e. the one PowerPoint's object model exposes through
# SlideMaster.Theme. The authored commit swaps the two parts' contents so
# the deck's effective/active theme becomes the plain "Office" palette
# (while the unused "Integral" palette ends up parked in theme1.xml,
# still only linked from the notes master).
#
# Reproduce that effect on the theme actually reachable through the
# object model by rewriting its twelve scheme colours (and, for
# completeness, its names) to the "Office" theme's values, in the
# ThemeColorScheme order PowerPoint uses:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink

function Convert-HexToOleColor([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$officeThemeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme

# Best-effort: line up the display names with the Office theme too.
try { $theme.Name = "Office Theme" } catch { }
try { $theme.ThemeColorScheme.Name = "Office" } catch { }
try { $theme.ThemeFontScheme.Name = "Office" } catch { }

$colorScheme = $theme.ThemeColorScheme
for ($i = 0; $i -lt $officeThemeColors.Length; $i++) {
    $colorScheme.Colors($i + 1).RGB = Convert-HexToOleColor $officeThemeColors[$i]
}
